$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("关于准予撤销上海市嘉定区朱家桥邮政支局的公告", "2025-12-05", "https://sh.spb.gov.cn/shsyzglj/c100057/c100058/202512/b487344d328e4d2fa163d4fe9a0fe502.shtml"),
    @("关于准予撤销上海市徐汇区柳州路邮政所的公告", "2025-12-02", "https://sh.spb.gov.cn/shsyzglj/c100057/c100058/202512/be58981880de42c7822366e7faabd2cb.shtml"),
    @("关于准予撤销上海市松江区达丰邮政所的公告", "2025-10-31", "https://sh.spb.gov.cn/shsyzglj/c100057/c100058/202510/e5a44a5099d1476fa0e479b321267ac3.shtml")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Style = "Normal"
    $row++
}
